# Finish section on Nijssen 07 Heuristic
# - adds a second (empty) worksheet "Tabelle2" after "Tabelle1"
# - appends a second example board (rows 12-20, columns D:L) to "Tabelle1",
#   mirroring the first example (rows 1-9) but with the corner-adjacent
#   squares (E/L on the inner rows, G/H/I/J on the border rows) filled in.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- add the new, empty "Tabelle2" sheet right after "Tabelle1" -----------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle2"

# --- duplicate the first example table (D1:L9) down to D12:L20 ------------
$src = $ws1.Range("D1:L9")
$dst = $ws1.Range("D12:L20")
$src.Copy($dst)

# row 1 has no D1 cell (it starts at E1), so the copy should not leave a
# stray blank D12 cell behind either
$ws1.Range("D12").ClearContents()

# --- fill in the previously-blank cells with the missing letters ----------
$ws1.Range("G13").Value = "A"
$ws1.Range("H13").Value = "B"
$ws1.Range("I13").Value = "B"
$ws1.Range("J13").Value = "A"

$ws1.Range("E15").Value = "A"
$ws1.Range("L15").Value = "A"

$ws1.Range("E16").Value = "B"
$ws1.Range("L16").Value = "B"

$ws1.Range("E17").Value = "B"
$ws1.Range("L17").Value = "B"

$ws1.Range("E18").Value = "A"
$ws1.Range("L18").Value = "A"

$ws1.Range("G20").Value = "A"
$ws1.Range("H20").Value = "B"
$ws1.Range("I20").Value = "B"
$ws1.Range("J20").Value = "A"

# --- restore view state: Tabelle1 active, scrolled to the new table -------
$ws1.Activate()
$ws1.Range("S20").Select()
